$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete row 17 (old Carmen row); row 18 (Samuel) shifts up to become new row 17,
# carrying its "last row" border styling with it.
$ws.Rows.Item(17).Delete()

# Re-enter Carmen's data into the new row 17 (replacing what shifted up),
# since Carmen remains in the updated dataset while Samuel is removed.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45592993"
$ws.Range("D17").Value = "CARMEN CECILIA RODRIGUEZ GONZALEZ"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1423500

# Update summary totals: one fewer worker and one fewer period now in the table,
# and the aggregate "Valor Mora" no longer includes Samuel's period value.
$ws.Range("E11").Value = 116000
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
